$wb = $excel.ActiveWorkbook

# --- Update the status text everywhere it occurs: "Ready for handoff" -> "In Translation" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value2 = "In Translation"
$wsOverview.Range("F2").Value2 = "In Translation"

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value2 = "In Translation"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value2 = "In Translation"

# --- Narrow the (now shorter) status columns to match the regenerated report's autofit widths ---
$wsOverview.Range("E1").ColumnWidth = 12.5
$wsOverview.Range("F1").ColumnWidth = 12.5

$wsZh.Range("C1").ColumnWidth = 12.5

$wsDe.Range("C1").ColumnWidth = 12.5
